$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 20, shifting rows 20..68 down to 21..69
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new weekly record
$ws.Range("A20").Value = 11
$ws.Range("B20").Value = "Vega Monumental Concepción"
$ws.Range("C20").Value = "Bíobío"
$ws.Range("D20").Value = 44497
$ws.Range("E20").Value = 8
$ws.Range("F20").Value = 100112032
$ws.Range("G20").Value = "Zapallo italiano"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 430
$ws.Range("K20").Value = 9500
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = 9767
$ws.Range("N20").Value = "$/caja 60 unidades"
$ws.Range("O20").Value = "Región de O'Higgins"
$ws.Range("P20").Value = 163
$ws.Range("Q20").Value = 60
$ws.Range("R20").Value = "Hortaliza"
